$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update announcement titles/summaries from the "200120" batch to "230120",
# tightening up spacing around the parenthetical group name, per the
# "member persona's regression" content refresh.
$ws.Range("A2").Value = "Announcement1-230120(UI Network (Spanish)"
$ws.Range("B2").Value = "Announcement1-230120Summary"

$ws.Range("A3").Value = "Announcement2-230120(Private Group)"
$ws.Range("B3").Value = "Announcement2-230120Summary"

$ws.Range("A4").Value = "Announcement3-230120(External Members)"
$ws.Range("B4").Value = "Announcement3-230120Summary"

$ws.Range("A5").Value = "Announcement4-230120(Network)"
$ws.Range("B5").Value = "Announcement4-230120Summary"

$ws.Range("A6").Value = "Announcement5-230120(Distributed Resiliency)"
$ws.Range("B6").Value = "Announcement5-230120Summary"

$ws.Range("A7").Value = "Announcement6-230120(Internal Admins)"
$ws.Range("B7").Value = "Announcement6-230120Summary"

# Widen column A to fit the updated announcement text.
$ws.Columns.Item(1).ColumnWidth = 42.62
